$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($col = 1; $col -le 11; $col++) {
    $c = $ws.Cells.Item(82, $col)
    Write-Host "COL=$col|Interior=$($c.Interior.Color)"
}
